$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day's data row (row 87) — daily update at 8 AM UTC
$ws.Range("A87").Value = 45827
$ws.Range("B87").Value = 367
$ws.Range("C87").Value = 377
$ws.Range("D87").Value = 376

# A86 loses its date-only format (style 3) in favor of the running
# date+time format (style 2); A87 takes on the date-only format instead.
$ws.Range("A86").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A87").NumberFormat = "YYYY-MM-DD"
